$d = $word.ActiveDocument

# The document ends with a paragraph that only contains "not finished yet".
# The edit appends a new run " :/" right after it, using the same
# Arial / 24-half-point (12pt) formatting as the existing text, but as its
# own separate run (matching how Word emits it) rather than merging into
# the existing run's text.

$lastPara = $d.Paragraphs.Last
$r = $lastPara.Range
$r.Collapse(0)              # wdCollapseEnd - move to the very end of the paragraph
$r.InsertAfter(" :/")       # $r now spans the newly inserted " :/" text

# Apply the same character formatting as the rest of the paragraph. We
# nudge the size away and back so the new text is recorded as a distinct
# run rather than being silently coalesced into the preceding run.
$r.Font.NameAscii = "Arial"
$r.Font.NameOther = "Arial"
$r.Font.NameBi = "Arial"
$r.Font.Size = 14
$r.Font.Size = 12
